$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 529.3077
$ws.Range("I9").Value = 377.375
$ws.Range("J9").Value = 772.4
$ws.Range("K9").Value = 377.375
$ws.Range("L9").Value = 772.4
$ws.Range("M9").Value = -208.375
$ws.Range("N9").Value = -1110.4
$ws.Range("H17").Value = 585.95557
$ws.Range("J17").Value = 585.95557
$ws.Range("L17").Value = 1757.86671
$ws.Range("N17").Value = -2093.86671
$ws.Range("H62").Value = 5091.5
$ws.Range("I62").Value = 5138.75
$ws.Range("K62").Value = 5138.75
$ws.Range("M62").Value = -4514.75
$ws.Range("H65").Value = 5091.5
$ws.Range("I65").Value = 5138.75
$ws.Range("K65").Value = 25693.75
$ws.Range("M65").Value = -22573.75
$ws.Range("H107").Value = 1520.6
$ws.Range("I107").Value = 1149.25
$ws.Range("K107").Value = 1149.25
$ws.Range("M107").Value = 770.75
$ws.Range("H112").Value = 4998.4443
$ws.Range("J112").Value = 4569.4287
$ws.Range("L112").Value = 13708.2861
$ws.Range("N112").Value = -15924.2861
$ws.Range("H121").Value = 655
$ws.Range("J121").Value = 655
$ws.Range("L121").Value = 1965
$ws.Range("N121").Value = -5459
$ws.Range("H125").Value = 1799.5
$ws.Range("J125").Value = 1327.4546
$ws.Range("L125").Value = 11947.0914
$ws.Range("N125").Value = -16867.0914
$ws.Range("H131").Value = 3844.9546
$ws.Range("I131").Value = 1294.2632
$ws.Range("J131").Value = 19999.334
$ws.Range("K131").Value = 3882.7896
$ws.Range("L131").Value = 59998.00199999999
$ws.Range("M131").Value = 1157.2104
$ws.Range("N131").Value = -70078.00199999999
$ws.Range("H137").Value = 1775.36
$ws.Range("I137").Value = 1363.3636
$ws.Range("J137").Value = 2099.0715
$ws.Range("K137").Value = 4090.0908
$ws.Range("L137").Value = 6297.2145
$ws.Range("M137").Value = -1540.0908
$ws.Range("N137").Value = -11397.2145
$ws.Range("H138").Value = 7487.4375
$ws.Range("I138").Value = 7487.4375
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 22462.3125
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -17322.3125
$ws.Range("N138").ClearContents()
$ws.Range("H141").Value = 1669.9
$ws.Range("I141").Value = 1385.7142
$ws.Range("K141").Value = 4157.142599999999
$ws.Range("M141").Value = 1022.857400000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32629.457
$ws.Range("I32").Value = 33972.355
$ws.Range("K32").Value = 33972.355
$ws.Range("M32").Value = -33685.355
$ws.Range("H61").Value = 4991.7144
$ws.Range("I61").Value = 3110.875
$ws.Range("J61").Value = 7499.5
$ws.Range("K61").Value = 3110.875
$ws.Range("L61").Value = 7499.5
$ws.Range("M61").Value = -2898.875
$ws.Range("N61").Value = -7923.5
$ws.Range("H97").Value = 7439.737
$ws.Range("I97").Value = 10070.909
$ws.Range("K97").Value = 10070.909
$ws.Range("M97").Value = -9574.909
$ws.Range("H132").Value = 45959.566
$ws.Range("I132").Value = 45959.566
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 137878.698
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -135348.698
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 4991.7144
$ws.Range("I136").Value = 3110.875
$ws.Range("J136").Value = 7499.5
$ws.Range("K136").Value = 9332.625
$ws.Range("L136").Value = 22498.5
$ws.Range("M136").Value = -6782.625
$ws.Range("N136").Value = -27598.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3848533.5
$ws.Range("I20").Value = 6251239
$ws.Range("K20").Value = 6251239
$ws.Range("M20").Value = -6250992
$ws.Range("H80").Value = 240.63637
$ws.Range("J80").Value = 219.11111
$ws.Range("L80").Value = 219.11111
$ws.Range("N80").Value = -2215.11111
$ws.Range("H83").Value = 240.63637
$ws.Range("J83").Value = 219.11111
$ws.Range("L83").Value = 1095.55555
$ws.Range("N83").Value = -11079.55555
$ws.Range("H86").Value = 3206
$ws.Range("I86").Value = 1579.9286
$ws.Range("K86").Value = 1579.9286
$ws.Range("M86").Value = -456.9286
$ws.Range("H89").Value = 3206
$ws.Range("I89").Value = 1579.9286
$ws.Range("K89").Value = 7899.643
$ws.Range("M89").Value = -2283.643
$ws.Range("H134").Value = 2954.9355
$ws.Range("I134").Value = 2171.7307
$ws.Range("J134").Value = 7027.6
$ws.Range("K134").Value = 6515.1921
$ws.Range("L134").Value = 21082.8
$ws.Range("M134").Value = -3980.1921
$ws.Range("N134").Value = -26152.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 911.8333
$ws.Range("I22").Value = 551.1539
$ws.Range("K22").Value = 551.1539
$ws.Range("M22").Value = -201.1539
$ws.Range("H31").Value = 1170.7333
$ws.Range("I31").Value = 1046.2
$ws.Range("K31").Value = 1046.2
$ws.Range("M31").Value = -751.2
$ws.Range("H34").Value = 1170.7333
$ws.Range("I34").Value = 1046.2
$ws.Range("K34").Value = 1046.2
$ws.Range("M34").Value = -844.2
$ws.Range("H58").Value = 57960.89
$ws.Range("J58").Value = 2975.75
$ws.Range("L58").Value = 2975.75
$ws.Range("N58").Value = -3381.75
$ws.Range("H62").Value = 3616.1428
$ws.Range("J62").Value = 2870
$ws.Range("L62").Value = 2870
$ws.Range("N62").Value = -4118
$ws.Range("H65").Value = 3616.1428
$ws.Range("J65").Value = 2870
$ws.Range("L65").Value = 14350
$ws.Range("N65").Value = -20590
$ws.Range("H122").Value = 2168.923
$ws.Range("I122").Value = 2183
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6549
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4099
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 2502.2222
$ws.Range("I132").Value = 3064.4
$ws.Range("J132").Value = 1799.5
$ws.Range("K132").Value = 9193.200000000001
$ws.Range("L132").Value = 5398.5
$ws.Range("M132").Value = -6663.200000000001
$ws.Range("N132").Value = -10458.5
$ws.Range("H136").Value = 57960.89
$ws.Range("J136").Value = 2975.75
$ws.Range("L136").Value = 8927.25
$ws.Range("N136").Value = -14027.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 4912.25
$ws.Range("I87").Value = 4912.25
$ws.Range("K87").Value = 14736.75
$ws.Range("M87").Value = -13488.75
$ws.Range("H90").Value = 4912.25
$ws.Range("I90").Value = 4912.25
$ws.Range("K90").Value = 44210.25
$ws.Range("M90").Value = -37970.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 74151.42999999999
$ws.Range("I132").Value = 79546.62
$ws.Range("K132").Value = 238639.86
$ws.Range("M132").Value = -236109.86

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 2245.9443
$ws.Range("J7").Value = 9999
$ws.Range("K7").Value = 2245.9443
$ws.Range("L7").Value = 9999
$ws.Range("M7").Value = -2133.9443
$ws.Range("N7").Value = -10223
$ws.Range("H26").Value = 1500
$ws.Range("I26").Value = 1500
$ws.Range("K26").Value = 1500
$ws.Range("M26").Value = -1205
$ws.Range("H46").Value = 5232.467
$ws.Range("J46").Value = 2894.2327
$ws.Range("L46").Value = 2894.2327
$ws.Range("N46").Value = -3270.2327
$ws.Range("H82").Value = 2054.4285
$ws.Range("I82").Value = 1836.25
$ws.Range("J82").Value = 2345.3333
$ws.Range("K82").Value = 1836.25
$ws.Range("L82").Value = 2345.3333
$ws.Range("M82").Value = -1475.25
$ws.Range("N82").Value = -3067.3333
$ws.Range("H85").Value = 2054.4285
$ws.Range("I85").Value = 1836.25
$ws.Range("J85").Value = 2345.3333
$ws.Range("K85").Value = 1836.25
$ws.Range("L85").Value = 2345.3333
$ws.Range("M85").Value = -588.25
$ws.Range("N85").Value = -4841.3333
$ws.Range("H122").Value = 3687.4688
$ws.Range("I122").Value = 3090.8635
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9272.5905
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6822.5905
$ws.Range("N122").Value = -19900
$ws.Range("I126").Value = 2245.9443
$ws.Range("J126").Value = 9999
$ws.Range("K126").Value = 6737.8329
$ws.Range("L126").Value = 29997
$ws.Range("M126").Value = -4267.8329
$ws.Range("N126").Value = -34937
$ws.Range("H132").Value = 61546.953
$ws.Range("I132").Value = 77780.19
$ws.Range("K132").Value = 233340.57
$ws.Range("M132").Value = -230810.57

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 2656
$ws.Range("I17").Value = 2656
$ws.Range("K17").Value = 2656
$ws.Range("M17").Value = -2484
$ws.Range("H62").Value = 172000.5
$ws.Range("I62").Value = 5333.6665
$ws.Range("K62").Value = 5333.6665
$ws.Range("M62").Value = -4709.6665
$ws.Range("H65").Value = 172000.5
$ws.Range("I65").Value = 5333.6665
$ws.Range("K65").Value = 26668.3325
$ws.Range("M65").Value = -23548.3325
$ws.Range("H122").Value = 2602
$ws.Range("I122").Value = 1804.5
$ws.Range("K122").Value = 5413.5
$ws.Range("M122").Value = -2963.5
$ws.Range("H126").Value = 45366.707
$ws.Range("I126").Value = 53349.8
$ws.Range("K126").Value = 160049.4
$ws.Range("M126").Value = -157579.4

